$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 895.7049
$ws.Range("J17").Value = 895.6731
$ws.Range("L17").Value = 2687.0193
$ws.Range("N17").Value = -3023.0193

$ws.Range("H137").Value = 4221.524
$ws.Range("I137").Value = 4464
$ws.Range("J137").Value = 2766.6667
$ws.Range("K137").Value = 13392
$ws.Range("L137").Value = 8300.000100000001
$ws.Range("M137").Value = -10842
$ws.Range("N137").Value = -13400.0001

$ws.Range("H138").Value = 3538.9546
$ws.Range("I138").Value = 2735.7
$ws.Range("J138").Value = 4208.3335
$ws.Range("K138").Value = 8207.099999999999
$ws.Range("L138").Value = 12625.0005
$ws.Range("M138").Value = -3067.099999999999
$ws.Range("N138").Value = -22905.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4874
$ws.Range("I2").Value = 3430
$ws.Range("J2").Value = 10650
$ws.Range("K2").Value = 3430
$ws.Range("L2").Value = 10650
$ws.Range("M2").Value = -3317
$ws.Range("N2").Value = -10876

$ws.Range("H32").Value = 18014.863
$ws.Range("I32").Value = 18733.418
$ws.Range("J32").Value = 9991
$ws.Range("K32").Value = 18733.418
$ws.Range("L32").Value = 9991
$ws.Range("M32").Value = -18446.418
$ws.Range("N32").Value = -10565

$ws.Range("H61").Value = 406356.44
$ws.Range("I61").Value = 345499.28
$ws.Range("J61").Value = 503727.9
$ws.Range("K61").Value = 345499.28
$ws.Range("L61").Value = 503727.9
$ws.Range("M61").Value = -345287.28
$ws.Range("N61").Value = -504151.9

$ws.Range("H74").Value = 129324.24
$ws.Range("I74").Value = 157086.56
$ws.Range("J74").Value = 48561.137
$ws.Range("K74").Value = 157086.56
$ws.Range("L74").Value = 48561.137
$ws.Range("M74").Value = -156212.56
$ws.Range("N74").Value = -50309.137

$ws.Range("H77").Value = 129324.24
$ws.Range("I77").Value = 157086.56
$ws.Range("J77").Value = 48561.137
$ws.Range("K77").Value = 785432.8
$ws.Range("L77").Value = 242805.685
$ws.Range("M77").Value = -781064.8
$ws.Range("N77").Value = -251541.685

$ws.Range("H102").Value = 1662.5
$ws.Range("I102").Value = 1143.8572
$ws.Range("J102").Value = 3477.75
$ws.Range("K102").Value = 1143.8572
$ws.Range("L102").Value = 3477.75
$ws.Range("M102").Value = 478.1428000000001
$ws.Range("N102").Value = -6721.75

$ws.Range("H116").Value = 4874
$ws.Range("I116").Value = 3430
$ws.Range("J116").Value = 10650
$ws.Range("K116").Value = 3430
$ws.Range("L116").Value = 10650
$ws.Range("M116").Value = -1136
$ws.Range("N116").Value = -15238

$ws.Range("H132").Value = 22865.04
$ws.Range("I132").Value = 30675.973
$ws.Range("J132").Value = 4118.8
$ws.Range("K132").Value = 92027.91900000001
$ws.Range("L132").Value = 12356.4
$ws.Range("M132").Value = -89497.91900000001
$ws.Range("N132").Value = -17416.4

$ws.Range("H136").Value = 406356.44
$ws.Range("I136").Value = 345499.28
$ws.Range("J136").Value = 503727.9
$ws.Range("K136").Value = 1036497.84
$ws.Range("L136").Value = 1511183.7
$ws.Range("M136").Value = -1033947.84
$ws.Range("N136").Value = -1516283.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4874
$ws.Range("I3").Value = 3430
$ws.Range("J3").Value = 10650
$ws.Range("K3").Value = 3430
$ws.Range("L3").Value = 10650
$ws.Range("M3").Value = -3316
$ws.Range("N3").Value = -10878

$ws.Range("H103").Value = 15920
$ws.Range("J103").Value = 15920
$ws.Range("L103").Value = 15920
$ws.Range("N103").Value = -18264

$ws.Range("H128").Value = 550
$ws.Range("I128").Value = 550
$ws.Range("K128").Value = 1650
$ws.Range("M128").Value = 840

$ws.Range("H134").Value = 3052.3696
$ws.Range("I134").Value = 2781.6365
$ws.Range("J134").Value = 3739.6155
$ws.Range("K134").Value = 8344.9095
$ws.Range("L134").Value = 11218.8465
$ws.Range("M134").Value = -5809.9095
$ws.Range("N134").Value = -16288.8465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2462.2922
$ws.Range("I31").Value = 1943.6904
$ws.Range("J31").Value = 3409.3044
$ws.Range("K31").Value = 1943.6904
$ws.Range("L31").Value = 3409.3044
$ws.Range("M31").Value = -1648.6904
$ws.Range("N31").Value = -3999.3044

$ws.Range("H34").Value = 2462.2922
$ws.Range("I34").Value = 1943.6904
$ws.Range("J34").Value = 3409.3044
$ws.Range("K34").Value = 1943.6904
$ws.Range("L34").Value = 3409.3044
$ws.Range("M34").Value = -1741.6904
$ws.Range("N34").Value = -3813.3044

$ws.Range("H132").Value = 2457.2666
$ws.Range("I132").Value = 1130.6
$ws.Range("J132").Value = 3783.9333
$ws.Range("K132").Value = 3391.8
$ws.Range("L132").Value = 11351.7999
$ws.Range("M132").Value = -861.7999999999997
$ws.Range("N132").Value = -16411.7999

$ws.Range("H134").Value = 1756.7693
$ws.Range("I134").Value = 1162
$ws.Range("K134").Value = 3486
$ws.Range("M134").Value = -951

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1991.6538
$ws.Range("I131").Value = 1818.5714
$ws.Range("J131").Value = 2055.4211
$ws.Range("K131").Value = 5455.7142
$ws.Range("L131").Value = 6166.263300000001
$ws.Range("M131").Value = -415.7142000000003
$ws.Range("N131").Value = -16246.2633

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6927.143
$ws.Range("I102").Value = 4229
$ws.Range("J102").Value = 13672.5
$ws.Range("K102").Value = 4229
$ws.Range("L102").Value = 13672.5
$ws.Range("M102").Value = -2607
$ws.Range("N102").Value = -16916.5

$ws.Range("H122").Value = 377202
$ws.Range("J122").Value = 2350
$ws.Range("L122").Value = 7050
$ws.Range("N122").Value = -11950

$ws.Range("H126").Value = 1368.4062
$ws.Range("I126").Value = 1136.5769
$ws.Range("J126").Value = 2373
$ws.Range("K126").Value = 3409.7307
$ws.Range("L126").Value = 7119
$ws.Range("M126").Value = -939.7307000000001
$ws.Range("N126").Value = -12059

$ws.Range("H132").Value = 4770.7334
$ws.Range("I132").Value = 6564.75
$ws.Range("J132").Value = 3574.7222
$ws.Range("K132").Value = 19694.25
$ws.Range("L132").Value = 10724.1666
$ws.Range("M132").Value = -17164.25
$ws.Range("N132").Value = -15784.1666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2807.0715
$ws.Range("I40").Value = 2366.6667
$ws.Range("K40").Value = 2366.6667
$ws.Range("M40").Value = -2230.6667

$ws.Range("H107").Value = 990
$ws.Range("I107").Value = 990
$ws.Range("K107").Value = 990
$ws.Range("M107").Value = 930

$ws.Range("H132").Value = 11681.174
$ws.Range("I132").Value = 4174.75
$ws.Range("J132").Value = 15684.6
$ws.Range("K132").Value = 12524.25
$ws.Range("L132").Value = 47053.8
$ws.Range("M132").Value = -9994.25
$ws.Range("N132").Value = -52113.8

$ws.Range("H136").Value = 3965.1133
$ws.Range("I136").Value = 2501.6365
$ws.Range("J136").Value = 6379.85
$ws.Range("K136").Value = 7504.9095
$ws.Range("L136").Value = 19139.55
$ws.Range("M136").Value = -4954.9095
$ws.Range("N136").Value = -24239.55

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1198.0834
$ws.Range("I122").Value = 1133.375
$ws.Range("K122").Value = 3400.125
$ws.Range("M122").Value = -950.125

$ws.Range("H132").Value = 2689.9644
$ws.Range("I132").Value = 1643.4667
$ws.Range("K132").Value = 4930.4001
$ws.Range("M132").Value = -2400.4001

$ws.Range("H136").Value = 17137528
$ws.Range("I136").Value = 25026584
$ws.Range("J136").Value = 528991.0600000001
$ws.Range("K136").Value = 75079752
$ws.Range("L136").Value = 1586973.18
$ws.Range("M136").Value = -75077202
$ws.Range("N136").Value = -1592073.18
